$wb = $excel.ActiveWorkbook

# zh-cn sheet: update "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 18:53:42"
$wsZhCn.Range("H2").Value = "2016-03-17 18:53:59"

# de-de sheet: update the same two columns.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 18:53:45"
$wsDeDe.Range("H2").Value = "2016-03-17 18:54:07"
